$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v0 = @'
[Jianlei%Cao%NULL%0,                            Xiaorong%Hu%NULL%1,                            Wenlin%Cheng%NULL%6,                            Lei%Yu%NULL%6,                            Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%6,                            Qiang%Liu%liuqiang@irm-cams.ac.cn%7,                            Qiang%Liu%liuqiang@irm-cams.ac.cn%0]
'@
$ws.Range("E2").Value = $v0

$v1 = @'
[Analysis of myocardial injury in patients with COVID-19 and association between concomitant cardiovascular diseases and severity of COVID-19].
'@
$ws.Range("C3").Value = $v1

$v2 = @'
Objective: To evaluate the cardiovascular damage of patients with COVID-19, and determine the correlation of serum N-terminal pro B-type natriuretic peptide (NT-proBNP) and cardiac troponin-I (cTnI) with the severity of COVID-19, and the impact of concomitant cardiovascular disease on severity of COVID-19 was also evaluated.
 Methods: A cross-sectional study was designed on 150 consecutive patients with COVID-19 in the fever clinic of Tongji Hospital in Wuhan from January to February in 2020, including 126 mild cases and 24 cases in critical care.
 Both univariate and multivariate logistic regression were used to analyze the correlation of past medical history including hypertension, diabetes and coronary heart disease (CHD) , as well as the levels of serum NT-proBNP and cTnI to the disease severity of COVID-19 patients.
 Results: Age, hypersensitive C-reactive protein(hs-CRP) and serum creatinine levels of the patients were higher in critical care cases than in mild cases(all P<0.05).
 Prevalence of male, elevated NT-proBNP and cTnI, hypertension and coronary heart disease were significantly higher in critical cases care patients than in the mild cases(all P<0.05).
 Univariate logistic regression analysis showed that age, male, elevated NT-proBNP, elevated cTnI, elevated hs-CRP, elevated serum creatinine, hypertension, and CHD were significantly correlated with critical disease status(all P<0.05).
 Multivariate logistic regression analysis showed that elevated cTnI(OR=26.909, 95%CI 4.086-177.226, P=0.001) and CHD (OR=16.609, 95%CI 2.288-120.577, P=0.005) were the independent risk factors of critical disease status.
 Conclusions: COVID-19 can significantly affect the heart function and lead to myocardial injury.
 The past medical history of CHD and increased level of cTnI are two independent determinants of clinical disease status in patients with COVID-19.
'@
$ws.Range("D3").Value = $v2

$v3 = @'
[ C%Chen%elasticNoEmail%0,   J T% Yan%elasticNoEmail%0,   N% Zhou%elasticNoEmail%0,   J P% Zhao%elasticNoEmail%0,   D W% Wang%elasticNoEmail%0]
'@
$ws.Range("E3").Value = $v3

$v4 = @'
10.3760/cma.j.cn112148-20200225-00123
'@
$ws.Range("F3").Value = $v4

$v5 = @'
ELASTIC
'@
$ws.Range("G3").Value = $v5

$v6 = @'
2020-03-06
'@
$ws.Range("H3").Value = $v6

$v7 = @'
[Chaolin%Huang%NULL%0,                            Yeming%Wang%NULL%14,                            Xingwang%Li%NULL%6,                            Lili%Ren%NULL%5,                            Jianping%Zhao%NULL%12,                            Yi%Hu%NULL%11,                            Li%Zhang%NULL%11,                            Guohui%Fan%NULL%14,                            Jiuyang%Xu%NULL%14,                            Xiaoying%Gu%NULL%14,                            Zhenshun%Cheng%NULL%7,                            Ting%Yu%NULL%24,                            Jiaan%Xia%NULL%6,                            Yuan%Wei%NULL%18,                            Wenjuan%Wu%NULL%5,                            Xuelei%Xie%NULL%6,                            Wen%Yin%NULL%9,                            Hui%Li%NULL%15,                            Min%Liu%NULL%5,                            Yan%Xiao%NULL%5,                            Hong%Gao%NULL%5,                            Li%Guo%NULL%6,                            Jungang%Xie%NULL%6,                            Guangfa%Wang%NULL%5,                            Rongmeng%Jiang%NULL%5,                            Zhancheng%Gao%NULL%5,                            Qi%Jin%NULL%5,                            Jianwei%Wang%wangjw28@163.com%5,                            Bin%Cao%caobin_ben@163.com%13]
'@
$ws.Range("E4").Value = $v7

$v8 = @'
[Clinical characteristics of 30 medical workers infected with new coronavirus pneumonia].
'@
$ws.Range("C5").Value = $v8

$v9 = @'
Objective: To investigate the clinical characteristics of medical staff with novel coronavirus pneumonia(NCP).
 Methods: 30 patients infected with novel coronavirus referred to jianghan university hospital between January 11, 2020 and January 3, 2020 were studied.
 The data reviewed included those of clinical manifestations, laboratory investigation and Radiographic features.
 Results: The patients consisted of 10 men and 20 women, including 22 doctors and 8 nurses,aged 21~59 years(mean 35±8 years).They were divided to 26 common type and 4 severe cases, all of whom had close(within 1m) contact with patients infected of novel coronavirus pneumonia.
 The average contact times were 12 (7,16) and the average cumulative contact time was 2 (1.5,2.7) h.Clinical symptoms of these patients were fever in 23 patients (76.67%) , headache in 16 petients (53.33%) , fatigue or myalgia in 21patients (70%) , nausea, vomiting or diarrhea in 9 petients (30%) , cough in 25 petients (83.33%) , and dyspnea in 14 petients (46.67%) .Routine blood test revealed WBC <4.0×10(9)/L in 8 petients (26.67%) , (4-10) ×10(9)/L in 22 petients (73.33%) , and WBC>4.0×10(9)/L in 4 petients (13.33%) during the disease.Lymphocyte count <1.0×10(9)/L occurred in 12 petients (40%),abnormal liver function in 7 petients (23.33%) ,myocardial damage in 5 petients(16.67%), elevated D-dimer (>0.5mg/l) in 5 patients (16.67%).
 Compared with normal patients, the average exposure times, cumulative exposure time, BMI, Fever time, white blood cell count, liver enzyme, LDH, myoenzyme and D-dimer were significantly increased in severe patients, while the lymphocyte count and albumin levels in peripheral blood were significantly decreased.Chest CT mainly showed patchy shadows and interstitial changes.According to imaging examination, 11 patients (36.67%) showed Unilateral pneumonia and 19 patients (63.33%) showed bilateral pneumonia,4 patients (13.33%) showed bilateral multiple mottling and ground-glass opacity.Compared with the patients infected in the protected period, the proportion of severe infection and bilateral pneumonia were both increased in the patients infected in unprotected period.
 Conclusion: Medical staffs are at higher risk of infection.Infection rates are associated with contact time, the amount of suction virus.
 Severe patients had BMI increased, heating time prolonged , white blood cell count, lymphocyte count, D-dimer and albumin level significantly changed and were prone to be complicated with liver damage and myocardial damage.Strict protection measures is important to prevent infection for medical workers.
'@
$ws.Range("D5").Value = $v9

$v10 = @'
[ M%Liu%elasticNoEmail%0,   P% He%elasticNoEmail%0,   H G% Liu%elasticNoEmail%0,   X J% Wang%elasticNoEmail%0,   F J% Li%elasticNoEmail%0,   S% Chen%elasticNoEmail%0,   J% Lin%elasticNoEmail%0,   P% Chen%elasticNoEmail%0,   J H% Liu%elasticNoEmail%0,   C H% Li%elasticNoEmail%0]
'@
$ws.Range("E5").Value = $v10

$v11 = @'
10.3760/cma.j.issn.1001-0939.2020.0016
'@
$ws.Range("F5").Value = $v11

$v12 = @'
ELASTIC
'@
$ws.Range("G5").Value = $v12

$v13 = @'
2020-02-17
'@
$ws.Range("H5").Value = $v13

$v14 = @'
[Qiurong%Ruan%NULL%0,                            Kun%Yang%NULL%8,                            Kun%Yang%NULL%0,                            Wenxia%Wang%NULL%8,                            Wenxia%Wang%NULL%0,                            Lingyu%Jiang%NULL%8,                            Lingyu%Jiang%NULL%0,                            Jianxin%Song%songsingsjx@sina.com%8,                            Jianxin%Song%songsingsjx@sina.com%0]
'@
$ws.Range("E6").Value = $v14

$v15 = @'
[Carb\u00f3-Cisnero%Yacquelin%coreGivesNoEmail%0,                           Fern\u00e1ndez-Gonz\u00e1lez%Paula%coreGivesNoEmail%3,                           Hierrezuelo-Rojas%Naifi%coreGivesNoEmail%3,                           Subert-Salas%Lizandra%coreGivesNoEmail%3]
'@
$ws.Range("E7").Value = $v15

$v16 = @'
[Aguilar%Medina Jos\u00e9 M.%coreGivesNoEmail%0,                           Claro%Vald\u00e9s Ram\u00f3n%coreGivesNoEmail%2,                           Obreg\u00f3n%Santos Angel G.%coreGivesNoEmail%2,                           Proh\u00edas%Mart\u00ednez Juan%coreGivesNoEmail%2,                           Rodr\u00edguez%Blanco Suilbert%coreGivesNoEmail%2]
'@
$ws.Range("E8").Value = $v16

$v17 = @'
[Clinical characteristics and outcomes of 112 cardiovascular disease patients infected by 2019-nCoV].
'@
$ws.Range("C9").Value = $v17

$v18 = @'
Objective: To explore the clinical characteristics and prognosis of the new coronavirus 2019-nCoV patients combined with cardiovascular disease (CVD).
 Methods: A retrospective analysis was performed on 112 COVID-19 patients with CVD admitted to the western district of Union Hospital in Wuhan, from January 20, 2020 to February 15, 2020.
 They were divided into critical group (ICU, n=16) and general group (n=96) according to the severity of the disease and patients were followed up to the clinical endpoint.
 The observation indicators included total blood count, C-reactive protein (CRP), arterial blood gas analysis, myocardial injury markers, coagulation function, liver and kidney function, electrolyte, procalcitonin (PCT), B-type natriuretic peptide (BNP), blood lipid, pulmonary CT and pathogen detection.
 Results: Compared with the general group, the lymphocyte count (0.74×10(9) (0.34×10(9), 0.94×10(9))/L vs.
 0.99×10(9) (0.71×10(9), 1.29×10(9))/L, P=0.03) was extremely lower in the critical group, CRP (106.98 (81.57, 135.76) mg/L vs.
 34.34 (9.55,76.54) mg/L, P<0.001) and PCT (0.20 (0.15,0.48) μg/L vs.
 0.11 (0.06,0.20)μg/L, P<0.001) were significantly higher in the critical group.
 The BMI of the critical group was significantly higher than that of the general group (25.5 (23.0, 27.5) kg/m(2) vs.
 22.0 (20.0, 24.0) kg/m(2), P=0.003).
 Patients were further divided into non-survivor group (17, 15.18%) group and survivor group (95, 84.82%).
 Among the non-survivors, there were 88.24% (15/17) patients with BMI> 25 kg/m(2), which was significantly higher than that of survivors (18.95% (18/95), P<0.001).
 Compared with the survived patients, oxygenation index (130 (102, 415) vs.
 434 (410, 444), P<0.001) was significantly lower and lactic acid (1.70 (1.30, 3.00) mmol/L vs.
 1.20 (1.10, 1.60) mmol/L, P<0.001) was significantly higher in the non-survivors.
 There was no significant difference in the proportion of ACEI/ARB medication between the critical group and the general group or between non-survivors and survivors (all P>0.05).
 Conclusion: COVID-19 patients combined with CVD are associated with a higher risk of mortality.
 Critical patients are characterized with lower lymphocyte counts.
 Higher BMI are more often seen in critical patients and non-survivor.
 ACEI/ARB use does not affect the morbidity and mortality of COVID-19 combined with CVD.
 Aggravating causes of death include fulminant inflammation, lactic acid accumulation and thrombotic events.
'@
$ws.Range("D9").Value = $v18

$v19 = @'
[ Y D%Peng%elasticNoEmail%0,   K% Meng%elasticNoEmail%0,   H Q% Guan%elasticNoEmail%0,   L% Leng%elasticNoEmail%0,   R R% Zhu%elasticNoEmail%0,   B Y% Wang%elasticNoEmail%0,   M A% He%elasticNoEmail%0,   L X% Cheng%elasticNoEmail%0,   K% Huang%elasticNoEmail%0,   Q T% Zeng%elasticNoEmail%0]
'@
$ws.Range("E9").Value = $v19

$v20 = @'
10.3760/cma.j.cn112148-20200220-00105
'@
$ws.Range("F9").Value = $v20

$v21 = @'
ELASTIC
'@
$ws.Range("G9").Value = $v21

$v22 = @'
2020-03-02
'@
$ws.Range("H9").Value = $v22

$v23 = @'
[Xiaobo%Yang%NULL%0,                            Yuan%Yu%NULL%8,                            Jiqian%Xu%NULL%8,                            Huaqing%Shu%NULL%6,                            Jia'an%Xia%NULL%10,                            Hong%Liu%NULL%6,                            Yongran%Wu%NULL%8,                            Lu%Zhang%NULL%6,                            Zhui%Yu%NULL%6,                            Minghao%Fang%NULL%6,                            Ting%Yu%NULL%0,                            Yaxin%Wang%NULL%8,                            Shangwen%Pan%NULL%6,                            Xiaojing%Zou%NULL%6,                            Shiying%Yuan%NULL%6,                            You%Shang%NULL%8]
'@
$ws.Range("E10").Value = $v23

$v24 = @'
[Fei%Zhou%NULL%0,                            Ting%Yu%NULL%0,                            Ronghui%Du%NULL%8,                            Guohui%Fan%NULL%0,                            Ying%Liu%NULL%13,                            Zhibo%Liu%NULL%9,                            Jie%Xiang%NULL%15,                            Yeming%Wang%NULL%0,                            Bin%Song%NULL%9,                            Xiaoying%Gu%NULL%0,                            Lulu%Guan%NULL%8,                            Yuan%Wei%NULL%0,                            Hui%Li%NULL%0,                            Xudong%Wu%NULL%9,                            Jiuyang%Xu%NULL%0,                            Shengjin%Tu%NULL%8,                            Yi%Zhang%NULL%9,                            Hua%Chen%NULL%9,                            Bin%Cao%NULL%0]
'@
$ws.Range("E11").Value = $v24

